$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header shared-string edits (preserve rich text via Characters) ---
$ws.Range("A8").Characters(21, 2).Text = "22"
$ws.Range("C9").Characters(47, 9).Text = "6/1/2025"
$ws.Range("C9").Characters(27, 9).Text = "5/26/2025"

# --- Model cells used as style/value sources (never modified by this script) ---
# C23 = style13 text "0" (shared string 20); E23 = style13 text "***.*" (shared string 21)
# C39 = style14 (integer); K39 = style15 (percent-like)

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 150
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 35.714285714285
$ws.Range("L15").Value = 171.428571428571
$ws.Range("M15").Value = 850
$ws.Range("N15").Value = 46.153846153846
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 80
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 37.037037037037
$ws.Range("I16").Value = 149
$ws.Range("J16").Value = 152
$ws.Range("K16").Value = -1.973684210526
$ws.Range("L16").Value = -27.669902912621
$ws.Range("M16").Value = 144.262295081967
$ws.Range("N16").Value = -85.727969348659
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 18.181818181818
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 2.380952380952
$ws.Range("I17").Value = 238
$ws.Range("J17").Value = 231
$ws.Range("K17").Value = 3.030303030303
$ws.Range("L17").Value = 14.423076923076
$ws.Range("M17").Value = 186.746987951807
$ws.Range("N17").Value = -13.454545454545
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -25.925925925925
$ws.Range("I18").Value = 179
$ws.Range("J18").Value = 166
$ws.Range("K18").Value = 7.831325301204
$ws.Range("L18").Value = -7.253886010362
$ws.Range("M18").Value = 31.617647058823
$ws.Range("N18").Value = -83.974932855863
$ws.Range("C19").Value = 39
$ws.Range("D19").Value = 31
$ws.Range("E19").Value = 25.806451612903
$ws.Range("F19").Value = 151
$ws.Range("G19").Value = 158
$ws.Range("H19").Value = -4.430379746835
$ws.Range("I19").Value = 738
$ws.Range("J19").Value = 872
$ws.Range("K19").Value = -15.366972477064
$ws.Range("L19").Value = -28
$ws.Range("M19").Value = -22.233930453108
$ws.Range("N19").Value = -81.067213955874
$ws.Range("D20").Value = 2
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -56.521739130434
$ws.Range("N20").Value = -94.011976047904
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 57
$ws.Range("E21").Value = 14.035087719298
$ws.Range("F21").Value = 258
$ws.Range("G21").Value = 262
$ws.Range("H21").Value = -1.526717557251
$ws.Range("I21").Value = 1333
$ws.Range("J21").Value = 1460
$ws.Range("K21").Value = -8.698630136986
$ws.Range("L21").Value = -19.939939939939
$ws.Range("M21").Value = 7.5
$ws.Range("N21").Value = -79.555214723926
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -25
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = 35.714285714285
$ws.Range("I22").Value = 95
$ws.Range("J22").Value = 79
$ws.Range("K22").Value = 20.253164556962
$ws.Range("L22").Value = -5
$ws.Range("M22").Value = 75.925925925925
$ws.Range("C24").Value = 84
$ws.Range("D24").Value = 87
$ws.Range("E24").Value = -3.448275862068
$ws.Range("F24").Value = 317
$ws.Range("H24").Value = -9.943181818181
$ws.Range("I24").Value = 1686
$ws.Range("J24").Value = 1779
$ws.Range("K24").Value = -5.227655986509
$ws.Range("L24").Value = 7.047619047619
$ws.Range("M24").Value = -11.820083682008
$ws.Range("C25").Value = 63
$ws.Range("D25").Value = 80
$ws.Range("E25").Value = -21.25
$ws.Range("F25").Value = 271
$ws.Range("G25").Value = 312
$ws.Range("H25").Value = -13.141025641025
$ws.Range("I25").Value = 1464
$ws.Range("J25").Value = 1569
$ws.Range("K25").Value = -6.692160611854
$ws.Range("L25").Value = -1.214574898785
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 33
$ws.Range("E26").Value = -45.454545454545
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 87
$ws.Range("H26").Value = -4.597701149425
$ws.Range("I26").Value = 451
$ws.Range("J26").Value = 426
$ws.Range("K26").Value = 5.868544600938
$ws.Range("L26").Value = 10.81081081081
$ws.Range("M26").Value = 78.260869565217
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = 22.222222222222
$ws.Range("L27").Value = 69.230769230769
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 29
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = 81.25
$ws.Range("I28").Value = 100
$ws.Range("J28").Value = 86
$ws.Range("K28").Value = 16.279069767441
$ws.Range("L28").Value = 13.636363636363
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 100

# --- Numeric -> Text cells (copy full cell incl. style+shared-string from model cells) ---
$ws.Range("C23").Copy($ws.Range("D15"))
$ws.Range("E23").Copy($ws.Range("E15"))
$ws.Range("C23").Copy($ws.Range("D27"))
$ws.Range("E23").Copy($ws.Range("E27"))
$ws.Range("C23").Copy($ws.Range("D31"))
$ws.Range("E23").Copy($ws.Range("E31"))

# --- Text -> Numeric cells (paste format from model cell, then set numeric value) ---
$ws.Range("C39").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = 1
$ws.Range("K39").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("C39").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("K39").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100
$ws.Range("C39").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("J33").Value = 1
$ws.Range("K39").Copy()
$ws.Range("K33").PasteSpecial(-4122)
$ws.Range("K33").Value = 0

$excel.CutCopyMode = 0
Write-Host "Edits applied"